$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# D-column (Price) values are forced to remain text (matching the source
# workbook, where every Price cell is stored as a string) by briefly
# switching the cell to text format before the write, then clearing the
# formatting again so no stray style index is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.080.69'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.478.90'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.44%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.44'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.63'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.24%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.512'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.476.15'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.137'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.331'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.928.58'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.30'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.931.54'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000169'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.428.81'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -5.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.01'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.43'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -7.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '350.16'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.01'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.06%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.44'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.22'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -6.01%  '
$ws.Range('E26').Value = '  -2.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.14'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -5.70%  '
$ws.Range('E28').Value = '  -33.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.586.67'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0900'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '508.53'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.64'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -6.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.77'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -4.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.23'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.61'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.58%  '
$ws.Range('E37').Value = '  -9.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.67'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.21'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -5.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.33'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -6.59%  '
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.68'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.81'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.97%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.326'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -4.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.35'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.72%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.82'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.83'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.48%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.513'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -5.19%  '
$ws.Range('B49').Value = 'Filecoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.44'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -6.61%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0251'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -7.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0728'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.44%  '
